{"js": "// The template paragraph contains a Word complex field whose field code is\n// the M2Doc query:  m:'doc.html'.fromHTMLURI()\n// (wrapped with leading/trailing spaces as instrText runs, delimited by\n// fldChar \"begin\"/\"end\"). This script rewrites that field into plain text\n// runs using the \"{ ... }\" token syntax expected by the\n// TokenIteratorFieldRewriterSplit parser, i.e. the single field becomes\n// seven plain <w:t> runs: \"{\", \"m\", \":\", \"'\", \"doc.html\",\n// \"'.fromHTMLURI()\" and \"}\".\n\nconst fields = context.document.body.fields;\nfields.load(\"items\");\nawait context.sync();\n\nif (fields.items.length === 0) {\n  throw new Error(\"Expected at least one field in the document body.\");\n}\n\nfields.load(\"items/code\");\nawait context.sync();\n\n// Find the field holding the M2Doc query defensively (don't assume index 0).\nfunction findQueryField(coll) {\n  for (let i = 0; i < coll.items.length; i++) {\n    const code = coll.items[i].code || \"\";\n    if (code.indexOf(\"fromHTMLURI\") !== -1) {\n      return coll.items[i];\n    }\n  }\n  return coll.items.length > 0 ? coll.items[0] : null;\n}\n\nconst target = findQueryField(fields);\n\n// Locate the paragraph that owns the field (the field itself has no\n// `getRange` in this host, so walk the paragraphs and check each one's\n// own `fields` collection).\nconst parentParas = target.parentBody.paragraphs;\nparentParas.load(\"items\");\nawait context.sync();\n\nlet ownerParagraph = null;\nfor (let i = 0; i < parentParas.items.length; i++) {\n  const pFields = parentParas.items[i].fields;\n  pFields.load(\"items\");\n  await context.sync();\n  if (pFields.items.length > 0) {\n    ownerParagraph = parentParas.items[i];\n    break;\n  }\n}\nif (!ownerParagraph) {\n  throw new Error(\"Could not locate the paragraph that owns the field.\");\n}\n\n// Insert the replacement plain-text runs immediately before the field's\n// content; inserting at the paragraph start (rather than replacing the\n// whole paragraph) keeps the paragraph's own identity/properties intact.\nconst startRange = ownerParagraph.getRange(\"Start\");\n\nconst replacementOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>{</w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>m</w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>:</w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>\\'</w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>doc.html</w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>\\'.fromHTMLURI()</w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nstartRange.insertOoxml(replacementOoxml, \"Before\");\nawait context.sync();\n\n// Re-fetch the field (it still exists, now preceded by the new runs) and\n// delete it, removing the old fldChar begin / instrText* / fldChar end run\n// sequence entirely.\nconst freshFields = context.document.body.fields;\nfreshFields.load(\"items\");\nawait context.sync();\nfreshFields.load(\"items/code\");\nawait context.sync();\n\nconst fieldToDelete = findQueryField(freshFields);\nif (fieldToDelete) {\n  fieldToDelete.delete();\n  await context.sync();\n}\n", "ps1": "# The template paragraph contains a Word complex field whose field code is\n# the M2Doc query:  m:'doc.html'.fromHTMLURI()\n# (wrapped with leading/trailing spaces as instrText runs, delimited by\n# fldChar \"begin\"/\"end\"). This script rewrites that field into plain text\n# runs using the \"{ ... }\" token syntax expected by the\n# TokenIteratorFieldRewriterSplit parser, i.e. the single field becomes\n# seven plain <w:t> runs: \"{\", \"m\", \":\", \"'\", \"doc.html\",\n# \"'.fromHTMLURI()\" and \"}\".\n\n$d = $word.ActiveDocument\n\nif ($d.Fields.Count -eq 0) {\n    throw \"Expected at least one field in the document.\"\n}\n\n# Find the field holding the M2Doc query defensively (don't assume it is\n# the first field in the document).\n$targetField = $null\nforeach ($fld in $d.Fields) {\n    if ($fld.Code.Text -like \"*fromHTMLURI*\") {\n        $targetField = $fld\n        break\n    }\n}\nif ($targetField -eq $null) {\n    $targetField = $d.Fields.Item(1)\n}\n\n# Find the paragraph that owns the field.\n$ownerPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Fields.Count -gt 0) {\n        $ownerPara = $p\n        break\n    }\n}\nif ($ownerPara -eq $null) {\n    throw \"Could not locate the paragraph that owns the field.\"\n}\n\n# Insert the replacement plain-text runs immediately before the field's\n# content; inserting at the paragraph start (rather than replacing the\n# whole paragraph range) keeps the paragraph's own identity/properties\n# (pPr, paraId, rsid, ...) intact.\n$paraRange = $ownerPara.Range\n$startRange = $d.Range($paraRange.Start, $paraRange.Start)\n\n$ooxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n  '<pkg:xmlData>' + `\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n  '<w:body>' + `\n  '<w:p>' + `\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>{</w:t></w:r>' + `\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>m</w:t></w:r>' + `\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>:</w:t></w:r>' + `\n  \"<w:r><w:rPr><w:lang w:val=\"\"en-US\"\"/></w:rPr><w:t>'</w:t></w:r>\" + `\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>doc.html</w:t></w:r>' + `\n  \"<w:r><w:rPr><w:lang w:val=\"\"en-US\"\"/></w:rPr><w:t>'.fromHTMLURI()</w:t></w:r>\" + `\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">}</w:t></w:r>' + `\n  '</w:p>' + `\n  '</w:body>' + `\n  '</w:document>' + `\n  '</pkg:xmlData>' + `\n  '</pkg:part>' + `\n  '</pkg:package>'\n\n$startRange.InsertXML($ooxml)\n\n# Re-resolve the field (the document changed; the field is still present,\n# now preceded by the new runs) and delete it, removing the old\n# fldChar begin / instrText* / fldChar end run sequence entirely.\n$fieldToDelete = $null\nforeach ($fld in $d.Fields) {\n    if ($fld.Code.Text -like \"*fromHTMLURI*\") {\n        $fieldToDelete = $fld\n        break\n    }\n}\nif ($fieldToDelete -eq $null -and $d.Fields.Count -gt 0) {\n    $fieldToDelete = $d.Fields.Item(1)\n}\nif ($fieldToDelete -ne $null) {\n    $fieldToDelete.Delete() | Out-Null\n}\n"}
